$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1): "Call Name *" -> "Call / Distribution Name *" ----
$ws.Range("C1").Value = "Call / Distribution Name *"

# ---- Existing data rows 2-6: fund renamed, folio numbers shortened ----
$ws.Range("A2").Value = "SAAS Fund"
$ws.Range("B2").Value = 6
$ws.Range("A3").Value = "SAAS Fund"
$ws.Range("B3").Value = 7
$ws.Range("A4").Value = "SAAS Fund"
$ws.Range("B4").Value = 8
$ws.Range("A5").Value = "SAAS Fund"
$ws.Range("B5").Value = 9
$ws.Range("C5").Value = "Call 1"
$ws.Range("A6").Value = "SAAS Fund"
$ws.Range("B6").Value = 10
$ws.Range("C6").Value = "Call 1"

# Apply the Arial font now used for the Fund (A) and Folio No (B) columns.
# Cells are touched one at a time (rather than via a single multi-cell Range)
# so each column settles on its own clean style record.
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Font.Name = "Arial"
    $ws.Cells.Item($r, 1).Font.Size = 11
}
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 2).Font.Name = "Arial"
    $ws.Cells.Item($r, 2).Font.Size = 11
}

# ---- New rows 7-11: "Distribution 1" entries mirroring rows 2-6 but with
#      negative quantities and a later issue date ----
$newRows = @(
    @{ Row=7;  Folio=6;  Call="Distribution 1"; Unit="Series C"; Qty=-5;  Prem=10 },
    @{ Row=8;  Folio=7;  Call="Distribution 1"; Unit="Series C"; Qty=-10; Prem=10 },
    @{ Row=9;  Folio=8;  Call="Distribution 1"; Unit="Series A"; Qty=-15; Prem=0  },
    @{ Row=10; Folio=9;  Call="Distribution 1"; Unit="Series A"; Qty=-20; Prem=0  },
    @{ Row=11; Folio=10; Call="Distribution 1"; Unit="Series B"; Qty=-25; Prem=5  }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = "SAAS Fund"
    $ws.Cells.Item($row, 2).Value = $r.Folio
    $ws.Cells.Item($row, 3).Value = $r.Call
    $ws.Cells.Item($row, 4).Value = $r.Unit
    $ws.Cells.Item($row, 5).Value = $r.Qty
    $ws.Cells.Item($row, 6).Value = 100
    $ws.Cells.Item($row, 7).Value = $r.Prem

    # Copy the date formatting (numFmtId 14) from an existing date cell so the
    # new cell reuses the same style record instead of creating a new one.
    $ws.Range("H2").Copy()
    $ws.Cells.Item($row, 8).PasteSpecial(-4122)
    $ws.Cells.Item($row, 8).Value = 44927

    $ws.Cells.Item($row, 9).Value = "No"
    $ws.Cells.Item($row, 10).Value = "Upload"
}

# Same Arial font treatment as rows 2-6 for the new Fund / Folio No cells
for ($r = 7; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Font.Name = "Arial"
    $ws.Cells.Item($r, 1).Font.Size = 11
}
for ($r = 7; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Font.Name = "Arial"
    $ws.Cells.Item($r, 2).Font.Size = 11
}

$excel.CutCopyMode = $false

# ---- Column widths: let Excel recompute the best-fit widths for the new data ----
$ws.Columns("A:C").AutoFit()

# ---- Selection mirrors the end of the author's edit session ----
$ws.Range("H11").Select()
